$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999999989342561
$ws.Range("E2").Value = 0.9999999989342561

$ws.Range("D3").Value = 0.9999983366317112
$ws.Range("E3").Value = 0.9999983366317112

$ws.Range("D4").Value = 0.9999999915616385
$ws.Range("E4").Value = 0.9999999915616385

$ws.Range("D5").Value = 0.003767655426160727
$ws.Range("E5").Value = 0.003767655426160727

$ws.Range("D6").Value = 0.003015807603838609
$ws.Range("E6").Value = 0.003015807603838609

$ws.Range("D7").Value = 0.9999980177652275
$ws.Range("E7").Value = [double]"1.982234772546221E-06"

$ws.Range("D8").Value = 0.9999999999999922
$ws.Range("E8").Value = [double]"7.771561172376096E-15"

$ws.Range("D9").Value = 0.9999999999996998
$ws.Range("E9").Value = [double]"3.002043058586423E-13"

$ws.Range("D10").Value = 0.8956182203807842
$ws.Range("E10").Value = 0.1043817796192158

$ws.Range("D11").Value = 0.9999999999973119
$ws.Range("E11").Value = [double]"2.688071987222429E-12"
$ws.Range("F11").Value = 5.267377853393555
